{"js": "// Replace the whole document body with a single paragraph reading\n// \"Looking forward to 2021\" and leave a collapsed \"_GoBack\" bookmark\n// right after the text (mirrors what Word stamps on save after an edit).\n\nconst body = context.document.body;\n\n// Wipe all existing paragraphs/content...\nbody.clear();\n\n// ...and type the new text into the (now single, empty) paragraph.\nbody.insertText(\"Looking forward to 2021\", Word.InsertLocation.start);\nawait context.sync();\n\n// Drop a collapsed \"_GoBack\" bookmark at the very end of the body, i.e.\n// right after the text we just inserted.\nconst end = body.getRange(Word.RangeLocation.end);\nend.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Replace the document body with a single paragraph reading\n# \"Looking forward to 2021\" and leave a collapsed \"_GoBack\" bookmark\n# right after the text (mirrors what Word stamps on save after an edit).\n\n$d = $word.ActiveDocument\n$newText = \"Looking forward to 2021\"\n\n# Collapse the document down to a single paragraph: repeatedly delete the\n# first paragraph's range (text + its own paragraph mark) so everything\n# merges into the final paragraph that remains.\nwhile ($d.Paragraphs.Count -gt 1) {\n    $d.Paragraphs.Item(1).Range.Delete()\n}\n\n# Overwrite the text of the one remaining paragraph.\n$p = $d.Paragraphs.Item(1)\n$p.Range.Text = $newText\n\n# Work out the character offset right after the new text (and before the\n# paragraph mark).\n$p = $d.Paragraphs.Item(1)\n$endPos = $p.Range.Start + $newText.Length\n\n# Temporarily append a sentinel character after the text. A collapsed\n# range sitting exactly at a paragraph/run boundary gets normalized to\n# span the whole adjacent run when used with Bookmarks.Add, so we add a\n# throwaway character to push our target position into the middle of the\n# run (a safe, unambiguous spot for a zero-length bookmark).\n$tail = $d.Range($endPos, $endPos)\n$tail.InsertAfter(\"X\")\n\n# Drop the collapsed \"_GoBack\" bookmark right after the visible text.\n$bmRange = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# Remove the sentinel character again.\n$d.Range($endPos, $endPos + 1).Delete()\n"}
